# Update cryptos list prices/volume-change figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store the Price column values as text
# (matching the original inlineStr cells) instead of re-parsing them as numbers.
$ws.Range("D2").Value = "'27.860.09"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").Value = "'1.627.45"
$ws.Range("E3").Value = "  -0.69%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'210.57"
$ws.Range("E5").Value = "  -0.98%  "

$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("D8").Value = "'23.35"
$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("E10").Value = "  -0.20%  "

$ws.Range("D11").Value = "'0.0880"
$ws.Range("E11").Value = "  -0.37%  "

$ws.Range("D12").Value = "'1.857.91"
$ws.Range("E12").Value = "  -0.63%  "

$ws.Range("D13").Value = "'1.625.10"
$ws.Range("E13").Value = "  -0.84%  "

$ws.Range("E14").Value = "  -1.41%  "

$ws.Range("D15").Value = "'0.561"
$ws.Range("E15").Value = "  -1.91%  "

$ws.Range("D16").Value = "'65.33"
$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("D17").Value = "'27.845.08"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").Value = "'229.81"
$ws.Range("E18").Value = "  -1.13%  "

$ws.Range("D19").Value = "'7.65"
$ws.Range("E19").Value = "  +0.80%  "

$ws.Range("E20").Value = "  -0.16%  "

$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").Value = "'4.32"
$ws.Range("E22").Value = "  -1.08%  "

$ws.Range("D23").Value = "'10.08"
$ws.Range("E23").Value = "  -3.62%  "

$ws.Range("E24").Value = "  -2.84%  "

$ws.Range("D25").Value = "'154.43"
$ws.Range("E25").Value = "  +1.15%  "

$ws.Range("D26").Value = "'6.90"
$ws.Range("E26").Value = "  +0.20%  "

$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("D28").Value = "'15.52"
$ws.Range("E28").Value = "  -1.08%  "

$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("E30").Value = "  -1.03%  "

$ws.Range("D31").Value = "'0.0482"
$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("E32").Value = "  +1.86%  "

$ws.Range("E33").Value = "  -0.46%  "

$ws.Range("D34").Value = "'1.396.78"
$ws.Range("E34").Value = "  -0.61%  "

$ws.Range("E35").Value = "  +0.46%  "

$ws.Range("D36").Value = "'1.01"
$ws.Range("E36").Value = "  +9.44%  "

$ws.Range("E38").Value = "  +0.23%  "

$ws.Range("E39").Value = "  -0.87%  "

$ws.Range("D40").Value = "'0.850"
$ws.Range("E40").Value = "  -3.46%  "

$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("E42").Value = "  -2.21%  "

$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").Value = "'65.75"
$ws.Range("E44").Value = "  -2.38%  "

$ws.Range("E45").Value = "  -1.43%  "

$ws.Range("D46").Value = "'1.767.82"
$ws.Range("E46").Value = "  -0.48%  "

$ws.Range("E47").Value = "  -2.72%  "

$ws.Range("D48").Value = "'88.03"
$ws.Range("E48").Value = "  +0.37%  "

$ws.Range("E49").Value = "  +1.47%  "

$ws.Range("E50").Value = "  +5.11%  "

$ws.Range("E51").Value = "  -0.67%  "
